# Update NATMI LR-pair (Vegfa -> Nrp2) sheet with newly recomputed TPM-based
# expression values. The new "raw" inputs are the ligand average/total
# expression per sending cluster (G/H) and the receptor average/total
# expression per target cluster (M/N); every other changed column
# (I/J, O/P, Q/R, S/T) is a derived specificity / product recomputed from
# those raw inputs, exactly mirroring the NATMI pipeline's own math.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new raw values, keyed by cluster name ---------------------------------
$clusters = @("ECs", "FAPs", "MuSCs")

# Ligand average / total expression value, per SENDING cluster (col A)
$ligAvg = @{
    "ECs"   = 1.854221333333333
    "FAPs"  = 45.11545066666667
    "MuSCs" = 7.310771333333332
}
$ligTot = @{
    "ECs"   = 5.562664
    "FAPs"  = 135.346352
    "MuSCs" = 21.932314
}

# Receptor average / total expression value, per TARGET cluster (col D)
$recAvg = @{
    "ECs"   = 25.37147633333333
    "FAPs"  = 10.21969166666667
    "MuSCs" = 8.297426666666667
}
$recTot = @{
    "ECs"   = 76.114429
    "FAPs"  = 30.659075
    "MuSCs" = 24.89228
}

# --- derived sums for specificity normalisation -----------------------------
$sumLigAvg = 0.0
$sumLigTot = 0.0
foreach ($c in $clusters) {
    $sumLigAvg += $ligAvg[$c]
    $sumLigTot += $ligTot[$c]
}

$sumRecAvg = 0.0
$sumRecTot = 0.0
foreach ($c in $clusters) {
    $sumRecAvg += $recAvg[$c]
    $sumRecTot += $recTot[$c]
}

# Edge weights (avg/total) for every sending x target combination, needed to
# normalise the edge-level specificity columns (S/T) over all 9 pairs.
$edgeAvg = @{}
$edgeTot = @{}
foreach ($s in $clusters) {
    foreach ($t in $clusters) {
        $edgeAvg["$s|$t"] = $ligAvg[$s] * $recAvg[$t]
        $edgeTot["$s|$t"] = $ligTot[$s] * $recTot[$t]
    }
}
$sumEdgeAvg = 0.0
$sumEdgeTot = 0.0
foreach ($s in $clusters) {
    foreach ($t in $clusters) {
        $sumEdgeAvg += $edgeAvg["$s|$t"]
        $sumEdgeTot += $edgeTot["$s|$t"]
    }
}

# --- row -> (sending cluster, target cluster) map, matching columns A & D --
$rows = @{
    2  = @{ S = "ECs";   T = "ECs"   }
    3  = @{ S = "ECs";   T = "FAPs"  }
    4  = @{ S = "ECs";   T = "MuSCs" }
    5  = @{ S = "FAPs";  T = "ECs"   }
    6  = @{ S = "FAPs";  T = "FAPs"  }
    7  = @{ S = "FAPs";  T = "MuSCs" }
    8  = @{ S = "MuSCs"; T = "ECs"   }
    9  = @{ S = "MuSCs"; T = "FAPs"  }
    10 = @{ S = "MuSCs"; T = "MuSCs" }
}

foreach ($r in $rows.Keys) {
    $s = $rows[$r].S
    $t = $rows[$r].T

    $g = $ligAvg[$s]
    $h = $ligTot[$s]
    $m = $recAvg[$t]
    $n = $recTot[$t]

    $i = $g / $sumLigAvg
    $j = $h / $sumLigTot

    $o = $m / $sumRecAvg
    $p = $n / $sumRecTot

    $q = $edgeAvg["$s|$t"]
    $rr = $edgeTot["$s|$t"]

    $sSpec = $q / $sumEdgeAvg
    $tSpec = $rr / $sumEdgeTot

    $ws.Cells.Item($r, 7).Value  = $g       # G - Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $h       # H - Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $i       # I - Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value = $j       # J - Ligand derived specificity (total)

    $ws.Cells.Item($r, 13).Value = $m       # M - Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $n       # N - Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $o       # O - Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value = $p       # P - Receptor derived specificity (total)

    $ws.Cells.Item($r, 17).Value = $q       # Q - Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $rr      # R - Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $sSpec   # S - Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $tSpec   # T - Edge total expression derived specificity
}
